# Neo4j deck - slide 14 ("TextBox 14" shape) - add inline "// comment" runs
# after two of the code-sample lines, per the target diff.
$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(14)
$sh = $s.Shapes.Item(6)
$tr = $sh.TextFrame.TextRange

# ---------------------------------------------------------------------
# Paragraph 1: "neo4j$ match(x) RETURN x"
#   -> append two tabs to the existing run, then add
#      "// match " (green) + "all nodes" (green) as two new runs.
# ---------------------------------------------------------------------
$para1 = $tr.Paragraphs(1, 1)
$run2  = $para1.Runs(2, 1)
$run2.Text = "match(x)" + [char]0x00A0 + "RETURN" + [char]0x00A0 + "x" + [char]9 + [char]9

$para1b = $tr.Paragraphs(1, 1)
$commentStart = $para1b.Start + $para1b.Length - 1
$para1b.InsertAfter("// match ") | Out-Null
$run3 = $tr.Characters($commentStart, 9)
$run3.Font.Color.RGB = 699961

$para1c = $tr.Paragraphs(1, 1)
$nodesStart = $para1c.Start + $para1c.Length - 1
$para1c.InsertAfter("all nodes") | Out-Null
$run4 = $tr.Characters($nodesStart, 9)
$run4.Font.Color.RGB = 699961

# ---------------------------------------------------------------------
# Paragraph 3: "neo4j$ match(x:Person) RETURN x"
#   -> append two tabs to the final "x" run, then add
#      "// match with label :Person" (green) as a new run.
# ---------------------------------------------------------------------
$para3 = $tr.Paragraphs(3, 1)
$run4b = $para3.Runs(4, 1)
$run4b.Text = "x" + [char]9 + [char]9

$para3b = $tr.Paragraphs(3, 1)
$labelStart = $para3b.Start + $para3b.Length - 1
$para3b.InsertAfter("// match with label :Person") | Out-Null
$run5 = $tr.Characters($labelStart, 28)
$run5.Font.Color.RGB = 699961
